# Added gesture module to automation
# The previous (now stale) Pass/Fail run results recorded in the
# "Results" column no longer apply once the gesture module test runs
# are introduced, so clear them out ahead of the new automation runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Clear the stale Results (column J) values for all data rows.
$ws.Range("J2:J8").ClearContents()

# Leave the selection on A2, as happens after the edits are made.
$ws.Range("A2").Select()
